# Auto-generated Excel COM-interop script
# Updates market-price derived columns (H-N) across all profession sheets
# to match the latest scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6: H6=284.33334, J6=533.25, L6=1599.75, N6=-1823.75
$ws.Range("H6").Value = 284.33334
$ws.Range("J6").Value = 533.25
$ws.Range("L6").Value = 1599.75
$ws.Range("N6").Value = -1823.75
# Row 17: H17=2695.9412, J17=2695.9412, L17=8087.823600000001, N17=-8423.8236
$ws.Range("H17").Value = 2695.9412
$ws.Range("J17").Value = 2695.9412
$ws.Range("L17").Value = 8087.823600000001
$ws.Range("N17").Value = -8423.8236
# Row 33: H33=314.33334, I33=116.30769, K33=116.30769, M33=112.69231
$ws.Range("H33").Value = 314.33334
$ws.Range("I33").Value = 116.30769
$ws.Range("K33").Value = 116.30769
$ws.Range("M33").Value = 112.69231
# Row 112: H112=2170.4375, J112=2464.818, L112=7394.454000000001, N112=-9610.454000000002
$ws.Range("H112").Value = 2170.4375
$ws.Range("J112").Value = 2464.818
$ws.Range("L112").Value = 7394.454000000001
$ws.Range("N112").Value = -9610.454000000002
# Row 116: H116=5684, I116=3696.6, J116=8996.333000000001, K116=3696.6, L116=8996.333000000001, M116=-254.5999999999999, N116=-15880.333
$ws.Range("H116").Value = 5684
$ws.Range("I116").Value = 3696.6
$ws.Range("J116").Value = 8996.333000000001
$ws.Range("K116").Value = 3696.6
$ws.Range("L116").Value = 8996.333000000001
$ws.Range("M116").Value = -254.5999999999999
$ws.Range("N116").Value = -15880.333
# Row 138: H138=6051.962, I138=7149.4287, J138=5654.6035, K138=21448.2861, L138=16963.8105, M138=-16308.2861, N138=-27243.8105
$ws.Range("H138").Value = 6051.962
$ws.Range("I138").Value = 7149.4287
$ws.Range("J138").Value = 5654.6035
$ws.Range("K138").Value = 21448.2861
$ws.Range("L138").Value = 16963.8105
$ws.Range("M138").Value = -16308.2861
$ws.Range("N138").Value = -27243.8105

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2=1847.5, I2=2004, J2=1299.75, K2=2004, L2=1299.75, M2=-1891, N2=-1525.75
$ws.Range("H2").Value = 1847.5
$ws.Range("I2").Value = 2004
$ws.Range("J2").Value = 1299.75
$ws.Range("K2").Value = 2004
$ws.Range("L2").Value = 1299.75
$ws.Range("M2").Value = -1891
$ws.Range("N2").Value = -1525.75
# Row 6: H6=47862.332, I6=113849.555, K6=113849.555, M6=-113676.555
$ws.Range("H6").Value = 47862.332
$ws.Range("I6").Value = 113849.555
$ws.Range("K6").Value = 113849.555
$ws.Range("M6").Value = -113676.555
# Row 32: H32=22805.191, I32=11909.381, K32=11909.381, M32=-11622.381
$ws.Range("H32").Value = 22805.191
$ws.Range("I32").Value = 11909.381
$ws.Range("K32").Value = 11909.381
$ws.Range("M32").Value = -11622.381
# Row 61: H61=2610.3333, I61=2610.3333, K61=2610.3333, M61=-2398.3333
$ws.Range("H61").Value = 2610.3333
$ws.Range("I61").Value = 2610.3333
$ws.Range("K61").Value = 2610.3333
$ws.Range("M61").Value = -2398.3333
# Row 97: H97=514.2593000000001, J97=632.25, L97=632.25, N97=-1624.25
$ws.Range("H97").Value = 514.2593000000001
$ws.Range("J97").Value = 632.25
$ws.Range("L97").Value = 632.25
$ws.Range("N97").Value = -1624.25
# Row 116: H116=1847.5, I116=2004, J116=1299.75, K116=2004, L116=1299.75, M116=290, N116=-5887.75
$ws.Range("H116").Value = 1847.5
$ws.Range("I116").Value = 2004
$ws.Range("J116").Value = 1299.75
$ws.Range("K116").Value = 2004
$ws.Range("L116").Value = 1299.75
$ws.Range("M116").Value = 290
$ws.Range("N116").Value = -5887.75
# Row 136: H136=2610.3333, I136=2610.3333, K136=7830.999899999999, M136=-5280.999899999999
$ws.Range("H136").Value = 2610.3333
$ws.Range("I136").Value = 2610.3333
$ws.Range("K136").Value = 7830.999899999999
$ws.Range("M136").Value = -5280.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3=1847.5, I3=2004, J3=1299.75, K3=2004, L3=1299.75, M3=-1890, N3=-1527.75
$ws.Range("H3").Value = 1847.5
$ws.Range("I3").Value = 2004
$ws.Range("J3").Value = 1299.75
$ws.Range("K3").Value = 2004
$ws.Range("L3").Value = 1299.75
$ws.Range("M3").Value = -1890
$ws.Range("N3").Value = -1527.75
# Row 64: H64=1699.4445, I64=1466.3334, K64=1466.3334, M64=-1241.3334
$ws.Range("H64").Value = 1699.4445
$ws.Range("I64").Value = 1466.3334
$ws.Range("K64").Value = 1466.3334
$ws.Range("M64").Value = -1241.3334
# Row 67: H67=1699.4445, I67=1466.3334, K67=1466.3334, M67=-686.3334
$ws.Range("H67").Value = 1699.4445
$ws.Range("I67").Value = 1466.3334
$ws.Range("K67").Value = 1466.3334
$ws.Range("M67").Value = -686.3334
# Row 94: H94=1095.2778, I94=1075.9375, J94=1250, K94=1075.9375, L94=1250, M94=-624.9375, N94=-2152
$ws.Range("H94").Value = 1095.2778
$ws.Range("I94").Value = 1075.9375
$ws.Range("J94").Value = 1250
$ws.Range("K94").Value = 1075.9375
$ws.Range("L94").Value = 1250
$ws.Range("M94").Value = -624.9375
$ws.Range("N94").Value = -2152
# Row 99: H99=1573.875, I99=1479.1333, K99=1479.1333, M99=18.86670000000004
$ws.Range("H99").Value = 1573.875
$ws.Range("I99").Value = 1479.1333
$ws.Range("K99").Value = 1479.1333
$ws.Range("M99").Value = 18.86670000000004
# Row 105: H105=3275.7441, I105=2626.2666, J105=4774.5386, K105=2626.2666, L105=4774.5386, M105=-879.2665999999999, N105=-8268.5386
$ws.Range("H105").Value = 3275.7441
$ws.Range("I105").Value = 2626.2666
$ws.Range("J105").Value = 4774.5386
$ws.Range("K105").Value = 2626.2666
$ws.Range("L105").Value = 4774.5386
$ws.Range("M105").Value = -879.2665999999999
$ws.Range("N105").Value = -8268.5386

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 99: H99=11933.6, J99=13309.9, L99=13309.9, N99=-16305.9
$ws.Range("H99").Value = 11933.6
$ws.Range("J99").Value = 13309.9
$ws.Range("L99").Value = 13309.9
$ws.Range("N99").Value = -16305.9
# Row 124: H124=39092, J124=39092, L124=39092, N124=-44002
$ws.Range("H124").Value = 39092
$ws.Range("J124").Value = 39092
$ws.Range("L124").Value = 39092
$ws.Range("N124").Value = -44002
# Row 126: H126=11933.6, J126=13309.9, L126=39929.7, N126=-44869.7
$ws.Range("H126").Value = 11933.6
$ws.Range("J126").Value = 13309.9
$ws.Range("L126").Value = 39929.7
$ws.Range("N126").Value = -44869.7
# Row 134: H134=2820.1785, I134=2182.6316, K134=6547.8948, M134=-4012.8948
$ws.Range("H134").Value = 2820.1785
$ws.Range("I134").Value = 2182.6316
$ws.Range("K134").Value = 6547.8948
$ws.Range("M134").Value = -4012.8948

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 56: H56=11906.667, I56=11906.667, K56=11906.667, M56=-11376.667
$ws.Range("H56").Value = 11906.667
$ws.Range("I56").Value = 11906.667
$ws.Range("K56").Value = 11906.667
$ws.Range("M56").Value = -11376.667
# Row 112: H112=2900, I112=2500, J112=3300, K112=7500, L112=9900, M112=-6392, N112=-12116
$ws.Range("H112").Value = 2900
$ws.Range("I112").Value = 2500
$ws.Range("J112").Value = 3300
$ws.Range("K112").Value = 7500
$ws.Range("L112").Value = 9900
$ws.Range("M112").Value = -6392
$ws.Range("N112").Value = -12116
# Row 131: H131=1486.138, I131=1020.8333, J131=1607.5217, K131=3062.4999, L131=4822.5651, M131=1977.5001, N131=-14902.5651
$ws.Range("H131").Value = 1486.138
$ws.Range("I131").Value = 1020.8333
$ws.Range("J131").Value = 1607.5217
$ws.Range("K131").Value = 3062.4999
$ws.Range("L131").Value = 4822.5651
$ws.Range("M131").Value = 1977.5001
$ws.Range("N131").Value = -14902.5651
# Row 137: H137=7821.3335, I137=9000, J137=7232, K137=27000, L137=21696, M137=-21900, N137=-31896
$ws.Range("H137").Value = 7821.3335
$ws.Range("I137").Value = 9000
$ws.Range("J137").Value = 7232
$ws.Range("K137").Value = 27000
$ws.Range("L137").Value = 21696
$ws.Range("M137").Value = -21900
$ws.Range("N137").Value = -31896
# Row 140: H140=4316.0625, I140=3504.1428, K140=10512.4284, M140=-5332.428400000001
$ws.Range("H140").Value = 4316.0625
$ws.Range("I140").Value = 3504.1428
$ws.Range("K140").Value = 10512.4284
$ws.Range("M140").Value = -5332.428400000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2: H2=1662, I2=325, K2=325, M2=-212
$ws.Range("H2").Value = 1662
$ws.Range("I2").Value = 325
$ws.Range("K2").Value = 325
$ws.Range("M2").Value = -212
# Row 24: H24=40325.145, J24=40325.145, L24=40325.145, N24=-40671.145
$ws.Range("H24").Value = 40325.145
$ws.Range("J24").Value = 40325.145
$ws.Range("L24").Value = 40325.145
$ws.Range("N24").Value = -40671.145
# Row 97: H97=1863.9584, I97=1998.3158, K97=1998.3158, M97=-1502.3158
$ws.Range("H97").Value = 1863.9584
$ws.Range("I97").Value = 1998.3158
$ws.Range("K97").Value = 1998.3158
$ws.Range("M97").Value = -1502.3158
# Row 122: H122=852225.3, I122=131992.5, K122=395977.5, M122=-393527.5
$ws.Range("H122").Value = 852225.3
$ws.Range("I122").Value = 131992.5
$ws.Range("K122").Value = 395977.5
$ws.Range("M122").Value = -393527.5
# Row 126: H126=4996.727, J126=4998.5, L126=14995.5, N126=-19935.5
$ws.Range("H126").Value = 4996.727
$ws.Range("J126").Value = 4998.5
$ws.Range("L126").Value = 14995.5
$ws.Range("N126").Value = -19935.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22=500, I22=500, K22=500, M22=-205
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
# Row 27: H27=500, I27=500, K27=500, M27=-393
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
# Row 100: H100=2878.2222, I100=3129.1428, K100=3129.1428, M100=-2588.1428
$ws.Range("H100").Value = 2878.2222
$ws.Range("I100").Value = 3129.1428
$ws.Range("K100").Value = 3129.1428
$ws.Range("M100").Value = -2588.1428
# Row 122: H122=7723.3335, I122=6585, K122=19755, M122=-17305
$ws.Range("H122").Value = 7723.3335
$ws.Range("I122").Value = 6585
$ws.Range("K122").Value = 19755
$ws.Range("M122").Value = -17305
# Row 132: H132=4845.45, I132=3491.7273, K132=10475.1819, M132=-7945.1819
$ws.Range("H132").Value = 4845.45
$ws.Range("I132").Value = 3491.7273
$ws.Range("K132").Value = 10475.1819
$ws.Range("M132").Value = -7945.1819
# Row 136: H136=4550.75, I136=4550.75, K136=13652.25, M136=-11102.25
$ws.Range("H136").Value = 4550.75
$ws.Range("I136").Value = 4550.75
$ws.Range("K136").Value = 13652.25
$ws.Range("M136").Value = -11102.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126: H126=89371.914, I126=129684.125, K126=389052.375, M126=-386582.375
$ws.Range("H126").Value = 89371.914
$ws.Range("I126").Value = 129684.125
$ws.Range("K126").Value = 389052.375
$ws.Range("M126").Value = -386582.375
# Row 132: H132=2454.5, I132=1939.6, K132=5818.799999999999, M132=-3288.799999999999
$ws.Range("H132").Value = 2454.5
$ws.Range("I132").Value = 1939.6
$ws.Range("K132").Value = 5818.799999999999
$ws.Range("M132").Value = -3288.799999999999
# Row 133: H133=92998.336, J133=92998.336, L133=92998.336, N133=-103118.336
$ws.Range("H133").Value = 92998.336
$ws.Range("J133").Value = 92998.336
$ws.Range("L133").Value = 92998.336
$ws.Range("N133").Value = -103118.336
# Row 136: H136=93748.63, I136=1655.75, K136=4967.25, M136=-2417.25
$ws.Range("H136").Value = 93748.63
$ws.Range("I136").Value = 1655.75
$ws.Range("K136").Value = 4967.25
$ws.Range("M136").Value = -2417.25
